$p = $ppt.ActivePresentation
try {
  $p.ApplyTemplate("Office Theme")
  Write-Output "ok"
} catch {
  Write-Output "ERR: $_"
}
